$wb = $excel.ActiveWorkbook

# --- Sheet1: update the Departing/Returning date values ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("D2").Value = "26/08/2019"
$ws1.Range("E2").Value = "30/10/2019"
$ws1.Range("D3").Value = "20/08/2019"
$ws1.Range("E3").Value = "01/10/2019"

# --- Sheet2: update the Departing/Returning date values ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("D2").Value = "26/08/2019"
$ws2.Range("E2").Value = "30/11/2019"
$ws2.Range("D3").Value = "20/08/2019"
$ws2.Range("E3").Value = "01/11/2019"

# --- Update selections (Sheet2 first, Sheet1 last so Sheet1 ends up active) ---
$ws2.Range("E4").Select()
$ws1.Range("B14").Select()
